$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure column D updates are kept as text (avoid Excel auto-converting numeric-looking strings)
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D48").NumberFormat = "@"

$ws.Range("D2").Value = "27.949.27"
$ws.Range("E2").Value = "  -0.17%  "
$ws.Range("D3").Value = "1.633.94"
$ws.Range("E3").Value = "  -0.91%  "
$ws.Range("E4").Value = "  -0.05%  "
$ws.Range("D5").Value = "212.08"
$ws.Range("E5").Value = "  -0.78%  "
$ws.Range("E6").Value = "  -0.80%  "
$ws.Range("D7").Value = "1.00"
$ws.Range("E7").Value = "  -0.08%  "
$ws.Range("D8").Value = "23.29"
$ws.Range("E8").Value = "  -0.65%  "
$ws.Range("E9").Value = "  -2.83%  "
$ws.Range("D10").Value = "0.0614"
$ws.Range("E10").Value = "  -0.08%  "
$ws.Range("E11").Value = "  +1.06%  "
$ws.Range("D12").Value = "1.864.79"
$ws.Range("E12").Value = "  -0.84%  "
$ws.Range("D13").Value = "1.643.44"
$ws.Range("E13").Value = "  -0.29%  "
$ws.Range("E14").Value = "  -0.40%  "
$ws.Range("D15").Value = "0.567"
$ws.Range("E15").Value = "  +0.18%  "
$ws.Range("D16").Value = "65.32"
$ws.Range("E16").Value = "  -0.48%  "
$ws.Range("D17").Value = "27.948.46"
$ws.Range("D18").Value = "230.87"
$ws.Range("E18").Value = "  -0.97%  "
$ws.Range("E19").Value = "  -0.16%  "
$ws.Range("D20").Value = "7.54"
$ws.Range("E20").Value = "  -2.11%  "
$ws.Range("E21").Value = "  -0.07%  "
$ws.Range("E22").Value = "  -0.66%  "
$ws.Range("D23").Value = "10.36"
$ws.Range("E23").Value = "  -3.04%  "
$ws.Range("D24").Value = "2.07"
$ws.Range("E24").Value = "  -4.01%  "
$ws.Range("D25").Value = "154.85"
$ws.Range("E25").Value = "  +1.73%  "
$ws.Range("D26").Value = "6.98"
$ws.Range("E26").Value = "  +0.76%  "
$ws.Range("D27").Value = "15.66"
$ws.Range("E27").Value = "  -0.82%  "
$ws.Range("E28").Value = "  -0.52%  "
$ws.Range("E29").Value = "  -0.06%  "
$ws.Range("D30").Value = "1.19"
$ws.Range("E30").Value = "  -0.97%  "
$ws.Range("E31").Value = "  -0.18%  "
$ws.Range("D32").Value = "3.39"
$ws.Range("E32").Value = "  +1.08%  "
$ws.Range("D33").Value = "1.406.45"
$ws.Range("E33").Value = "  -2.73%  "
$ws.Range("E34").Value = "  -0.08%  "
$ws.Range("E35").Value = "  -0.01%  "
$ws.Range("E36").Value = "  +9.84%  "
$ws.Range("E37").Value = "  +1.41%  "
$ws.Range("E38").Value = "  +0.71%  "
$ws.Range("D39").Value = "0.561"
$ws.Range("E39").Value = "  +0.15%  "
$ws.Range("D40").Value = "0.873"
$ws.Range("E40").Value = "  -1.89%  "
$ws.Range("E41").Value = "  -0.11%  "
$ws.Range("E42").Value = "  -0.07%  "
$ws.Range("D43").Value = "66.96"
$ws.Range("E43").Value = "  -3.48%  "
$ws.Range("E44").Value = "  +2.68%  "
$ws.Range("D45").Value = "1.82"
$ws.Range("E45").Value = "  +1.00%  "
$ws.Range("E46").Value = "  -1.35%  "
$ws.Range("D47").Value = "1.774.76"
$ws.Range("E47").Value = "  -0.86%  "
$ws.Range("D48").Value = "87.85"
$ws.Range("E48").Value = "  -1.27%  "
$ws.Range("E49").Value = "  +0.78%  "
$ws.Range("E51").Value = "  -0.30%  "
